$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 216115.5
$ws.Range("I15").Value = 216115.5
$ws.Range("K15").Value = 648346.5
$ws.Range("M15").Value = -648177.5
$ws.Range("H18").Value = 1335.875
$ws.Range("I18").Value = 238
$ws.Range("J18").Value = 2433.75
$ws.Range("K18").Value = 238
$ws.Range("L18").Value = 2433.75
$ws.Range("M18").Value = 46
$ws.Range("N18").Value = -3001.75
$ws.Range("H40").Value = 2177.0715
$ws.Range("J40").Value = 1985.5714
$ws.Range("L40").Value = 1985.5714
$ws.Range("N40").Value = -2335.5714
$ws.Range("H98").Value = 321716.12
$ws.Range("I98").Value = 363134.53
$ws.Range("J98").Value = 723.5
$ws.Range("K98").Value = 363134.53
$ws.Range("L98").Value = 723.5
$ws.Range("M98").Value = -361636.53
$ws.Range("N98").Value = -3719.5
$ws.Range("H122").Value = 321716.12
$ws.Range("I122").Value = 363134.53
$ws.Range("J122").Value = 723.5
$ws.Range("K122").Value = 1089403.59
$ws.Range("L122").Value = 2170.5
$ws.Range("M122").Value = -1086953.59
$ws.Range("N122").Value = -7070.5
$ws.Range("H130").Value = 57516.332
$ws.Range("J130").Value = 57516.332
$ws.Range("L130").Value = 57516.332
$ws.Range("N130").Value = -67556.33199999999
$ws.Range("H137").Value = 23257100
$ws.Range("I137").Value = 31251038
$ws.Range("J137").Value = 2009.091
$ws.Range("K137").Value = 93753114
$ws.Range("L137").Value = 6027.272999999999
$ws.Range("M137").Value = -93750564
$ws.Range("N137").Value = -11127.273
$ws.Range("H138").Value = 8734517
$ws.Range("J138").Value = 10640991
$ws.Range("L138").Value = 31922973
$ws.Range("N138").Value = -31933253
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16346.902
$ws.Range("I32").Value = 3123.9595
$ws.Range("J32").Value = 70707.89
$ws.Range("K32").Value = 3123.9595
$ws.Range("L32").Value = 70707.89
$ws.Range("M32").Value = -2836.9595
$ws.Range("N32").Value = -71281.89
$ws.Range("H138").Value = 62429
$ws.Range("J138").Value = 62429
$ws.Range("L138").Value = 62429
$ws.Range("N138").Value = -72709
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 360508.25
$ws.Range("I105").Value = 3241.7646
$ws.Range("K105").Value = 3241.7646
$ws.Range("M105").Value = -1494.7646
$ws.Range("H134").Value = 4302.923
$ws.Range("I134").Value = 2811.8262
$ws.Range("J134").Value = 6446.375
$ws.Range("K134").Value = 8435.4786
$ws.Range("L134").Value = 19339.125
$ws.Range("M134").Value = -5900.4786
$ws.Range("N134").Value = -24409.125
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1398.6
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = 1248.25
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 1248.25
$ws.Range("M16").Value = -1713
$ws.Range("N16").Value = -1822.25
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = $null
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").Value = $null
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = $null
$ws.Range("H68").Value = 19165
$ws.Range("J68").Value = 19998
$ws.Range("L68").Value = 19998
$ws.Range("N68").Value = -21496
$ws.Range("H71").Value = 19165
$ws.Range("J71").Value = 19998
$ws.Range("L71").Value = 59994
$ws.Range("N71").Value = -67482
$ws.Range("H74").Value = 17796.25
$ws.Range("J74").Value = 17796.25
$ws.Range("L74").Value = 17796.25
$ws.Range("N74").Value = -19544.25
$ws.Range("H77").Value = 17796.25
$ws.Range("J77").Value = 17796.25
$ws.Range("L77").Value = 53388.75
$ws.Range("N77").Value = -62124.75
$ws.Range("H86").Value = 55558804
$ws.Range("I86").Value = 166668820
$ws.Range("J86").Value = 3797.8333
$ws.Range("K86").Value = 166668820
$ws.Range("L86").Value = 3797.8333
$ws.Range("M86").Value = -166667697
$ws.Range("N86").Value = -6043.8333
$ws.Range("H89").Value = 55558804
$ws.Range("I89").Value = 166668820
$ws.Range("J89").Value = 3797.8333
$ws.Range("K89").Value = 833344100
$ws.Range("L89").Value = 18989.1665
$ws.Range("M89").Value = -833338484
$ws.Range("N89").Value = -30221.1665
$ws.Range("H113").Value = 1398.6
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 1248.25
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 1248.25
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -5588.25
$ws.Range("H134").Value = 3113.8333
$ws.Range("I134").Value = 1584.9231
$ws.Range("J134").Value = 7089
$ws.Range("K134").Value = 4754.7693
$ws.Range("L134").Value = 21267
$ws.Range("M134").Value = -2219.7693
$ws.Range("N134").Value = -26337
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1084.5272
$ws.Range("I5").Value = 677.9286
$ws.Range("J5").Value = 1506.1852
$ws.Range("K5").Value = 2033.7858
$ws.Range("L5").Value = 4518.5556
$ws.Range("M5").Value = -1921.7858
$ws.Range("N5").Value = -4742.5556
$ws.Range("H54").Value = 4950
$ws.Range("J54").Value = 4950
$ws.Range("L54").Value = 14850
$ws.Range("N54").Value = -15968
$ws.Range("H109").Value = 2705.4
$ws.Range("I109").Value = 2013.5
$ws.Range("J109").Value = 3166.6667
$ws.Range("K109").Value = 6040.5
$ws.Range("L109").Value = 9500.000100000001
$ws.Range("M109").Value = -5000.5
$ws.Range("N109").Value = -11580.0001
$ws.Range("H131").Value = 7408942.5
$ws.Range("I131").Value = 554.2857
$ws.Range("J131").Value = 8773646
$ws.Range("K131").Value = 1662.8571
$ws.Range("L131").Value = 26320938
$ws.Range("M131").Value = 3377.1429
$ws.Range("N131").Value = -26331018
$ws.Range("H132").Value = 1087.6207
$ws.Range("I132").Value = 830.8889
$ws.Range("J132").Value = 1507.7273
$ws.Range("K132").Value = 7478.0001
$ws.Range("L132").Value = 13569.5457
$ws.Range("M132").Value = -4948.0001
$ws.Range("N132").Value = -18629.5457
$ws.Range("H135").Value = 1084.5272
$ws.Range("I135").Value = 677.9286
$ws.Range("J135").Value = 1506.1852
$ws.Range("K135").Value = 6101.3574
$ws.Range("L135").Value = 13555.6668
$ws.Range("M135").Value = -3566.3574
$ws.Range("N135").Value = -18625.6668
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I57").Value = 30000
$ws.Range("K57").Value = 30000
$ws.Range("M57").Value = -29180
$ws.Range("H70").Value = 4382.967
$ws.Range("I70").Value = 4341.72
$ws.Range("K70").Value = 4341.72
$ws.Range("M70").Value = -4071.72
$ws.Range("H73").Value = 4382.967
$ws.Range("I73").Value = 4341.72
$ws.Range("K73").Value = 4341.72
$ws.Range("M73").Value = -3405.72
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = $null
$ws.Range("N113").Value = -5340
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3284.4736
$ws.Range("I7").Value = 2400
$ws.Range("J7").Value = 3520.3333
$ws.Range("K7").Value = 2400
$ws.Range("L7").Value = 3520.3333
$ws.Range("M7").Value = -2288
$ws.Range("N7").Value = -3744.3333
$ws.Range("H92").Value = 30500
$ws.Range("J92").Value = 30500
$ws.Range("L92").Value = 30500
$ws.Range("N92").Value = -35492
$ws.Range("H122").Value = 3650
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3650
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 10950
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = -15850
$ws.Range("H126").Value = 3284.4736
$ws.Range("I126").Value = 2400
$ws.Range("J126").Value = 3520.3333
$ws.Range("K126").Value = 7200
$ws.Range("L126").Value = 10560.9999
$ws.Range("M126").Value = -4730
$ws.Range("N126").Value = -15500.9999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 16185.772
$ws.Range("I62").Value = 17835.273
$ws.Range("K62").Value = 17835.273
$ws.Range("M62").Value = -17211.273
$ws.Range("H65").Value = 16185.772
$ws.Range("I65").Value = 17835.273
$ws.Range("K65").Value = 89176.36500000001
$ws.Range("M65").Value = -86056.36500000001
$ws.Range("H81").Value = 3592.72
$ws.Range("J81").Value = 4621.8823
$ws.Range("L81").Value = 9243.7646
$ws.Range("N81").Value = -11365.7646
$ws.Range("H84").Value = 3592.72
$ws.Range("J84").Value = 4621.8823
$ws.Range("L84").Value = 46218.823
$ws.Range("N84").Value = -56826.823
$ws.Range("H86").Value = 45000
$ws.Range("J86").Value = 45000
$ws.Range("L86").Value = 45000
$ws.Range("N86").Value = -47246
$ws.Range("H89").Value = 45000
$ws.Range("J89").Value = 45000
$ws.Range("L89").Value = 225000
$ws.Range("N89").Value = -236232
$ws.Range("H136").Value = 1535.3704
$ws.Range("I136").Value = 777.25
$ws.Range("J136").Value = 3701.4285
$ws.Range("K136").Value = 2331.75
$ws.Range("L136").Value = 11104.2855
$ws.Range("M136").Value = 218.25
$ws.Range("N136").Value = -16204.2855
